$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A8: "Volume 30   Number  14" -> "...15" (last run text only) ---
$full = $ws.Range("A8").Text
$start = $full.IndexOf("14") + 1
$len = "14".Length
$ws.Range("A8").Characters($start, $len).Text = "15"

# --- C9: "Report Covering the Week  4/3/2023  Through  4/9/2023" -> dates +7 days ---
$full = $ws.Range("C9").Text
$start = $full.IndexOf("4/3/2023") + 1
$len = "4/3/2023".Length
$ws.Range("C9").Characters($start, $len).Text = "4/10/2023"

$full = $ws.Range("C9").Text
$start = $full.IndexOf("4/9/2023") + 1
$len = "4/9/2023".Length
$ws.Range("C9").Characters($start, $len).Text = "4/16/2023"

# --- Row 18: C18 changes from text "0" to number 4 (style 14 -> 15) ---
$ws.Range("C16").Copy($ws.Range("C18"))
$ws.Range("C18").Value = 4

# --- Row 27: C27 changes from text "0" to number 4 (style 14 -> 15) ---
$ws.Range("C16").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 4

# --- Row 30: C30 changes from number 1 to text "0" (style 15 -> 14, shared string) ---
$ws.Range("C14").Copy($ws.Range("C30"))

# --- Row 30: D30 changes from text "0" to number 1 (style 14 -> 15) ---
$ws.Range("D16").Copy($ws.Range("D30"))
$ws.Range("D30").Value = 1

# --- Row 30: E30 changes from text "***.* " to number -100 (style 14 -> 16) ---
$ws.Range("E16").Copy($ws.Range("E30"))
$ws.Range("E30").Value = -100

# --- Bulk numeric cell updates ---
$ws.Range("G15").Value = 1
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -50
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -50
$ws.Range("I16").Value = 38
$ws.Range("J16").Value = 53
$ws.Range("K16").Value = -28.301886792452
$ws.Range("L16").Value = 90
$ws.Range("M16").Value = 5.555555555555
$ws.Range("N16").Value = -89.645776566757
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 16.666666666666
$ws.Range("F17").Value = 19
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = 46.153846153846
$ws.Range("I17").Value = 57
$ws.Range("J17").Value = 50
$ws.Range("K17").Value = 14
$ws.Range("L17").Value = 39.024390243902
$ws.Range("M17").Value = 54.054054054054
$ws.Range("N17").Value = -58.695652173913
$ws.Range("D18").Value = 9
$ws.Range("E18").Value = -55.555555555555
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 26
$ws.Range("H18").Value = -50
$ws.Range("I18").Value = 42
$ws.Range("J18").Value = 67
$ws.Range("K18").Value = -37.31343283582
$ws.Range("L18").Value = -10.63829787234
$ws.Range("M18").Value = -27.586206896551
$ws.Range("N18").Value = -93.247588424437
$ws.Range("C19").Value = 28
$ws.Range("D19").Value = 38
$ws.Range("E19").Value = -26.315789473684
$ws.Range("F19").Value = 131
$ws.Range("G19").Value = 146
$ws.Range("H19").Value = -10.273972602739
$ws.Range("I19").Value = 515
$ws.Range("J19").Value = 492
$ws.Range("K19").Value = 4.674796747967
$ws.Range("L19").Value = 101.171875
$ws.Range("M19").Value = 17.045454545454
$ws.Range("N19").Value = -73.989898989899
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = -11.111111111111
$ws.Range("I20").Value = 21
$ws.Range("J20").Value = 28
$ws.Range("K20").Value = -25
$ws.Range("L20").Value = 40
$ws.Range("M20").Value = 90.90909090909
$ws.Range("N20").Value = -86
$ws.Range("C21").Value = 44
$ws.Range("D21").Value = 61
$ws.Range("E21").Value = -27.868852459016
$ws.Range("F21").Value = 179
$ws.Range("G21").Value = 211
$ws.Range("H21").Value = -15.165876777251
$ws.Range("I21").Value = 675
$ws.Range("J21").Value = 696
$ws.Range("K21").Value = -3.01724137931
$ws.Range("L21").Value = 75.78125
$ws.Range("M21").Value = 13.827993254637
$ws.Range("N21").Value = -79.383017715332
$ws.Range("C22").Value = 1
$ws.Range("E22").Value = -66.666666666666
$ws.Range("G22").Value = 9
$ws.Range("H22").Value = -66.666666666666
$ws.Range("I22").Value = 18
$ws.Range("J22").Value = 20
$ws.Range("K22").Value = -10
$ws.Range("L22").Value = 28.571428571428
$ws.Range("C24").Value = 51
$ws.Range("D24").Value = 40
$ws.Range("E24").Value = 27.5
$ws.Range("F24").Value = 185
$ws.Range("G24").Value = 153
$ws.Range("H24").Value = 20.915032679738
$ws.Range("I24").Value = 705
$ws.Range("J24").Value = 601
$ws.Range("K24").Value = 17.304492512479
$ws.Range("L24").Value = 70.702179176755
$ws.Range("M24").Value = 43.292682926829
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = 10
$ws.Range("F25").Value = 59
$ws.Range("H25").Value = 40.47619047619
$ws.Range("I25").Value = 185
$ws.Range("J25").Value = 137
$ws.Range("K25").Value = 35.036496350365
$ws.Range("L25").Value = 81.372549019607
$ws.Range("M25").Value = 28.472222222222
$ws.Range("J26").Value = 11
$ws.Range("K26").Value = -45.454545454545
$ws.Range("L26").Value = 50
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = 33.333333333333
$ws.Range("F27").Value = 7
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 22
$ws.Range("J27").Value = 28
$ws.Range("K27").Value = -21.428571428571
$ws.Range("L27").Value = 29.411764705882
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = -50
$ws.Range("J30").Value = 3
$ws.Range("K30").Value = 0
